$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("K7:L7").Copy($ws.Range("L7"))
for ($c = 10; $c -le 13; $c++) {
  $cell = $ws.Cells.Item(7, $c)
  Write-Host "col" $c "HAlign=" $cell.HorizontalAlignment "FontColor=" $cell.Font.Color
}
